$d = $word.ActiveDocument

# The paragraph currently reads (across several runs):
#   "Si alguien dice algo que mi hijo/a piensa que está equivocado e es malo,
#    mi hijo/s se siente asustado/a diciendo o que piensa si esa persona es:"
# It needs to become:
#   "Si alguien dice algo que mi hijo/a piensa que es equivocado o malo,
#    mi hijo/a se siente asustado/a de decir lo que piensa si esa persona es:"

$old = "Si alguien dice algo que mi hijo/a piensa que está equivocado e es malo, mi hijo/s se siente asustado/a diciendo o que piensa si esa persona es"
$new = "Si alguien dice algo que mi hijo/a piensa que es equivocado o malo, mi hijo/a se siente asustado/a de decir lo que piensa si esa persona es"

# Replace the body of the sentence, leaving the trailing colon untouched so
# it stays in its own run (matching the target markup).
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# Re-locate the replaced range and nudge its formatting (no-op toggle) so the
# run boundary is preserved between the new sentence text and the colon that
# follows it, instead of the two being auto-merged into a single run.
$r = $d.Content
$r.Find.Execute($new) | Out-Null
$r.Font.Bold = 1
$r.Font.Bold = 0
